# Auto-generated edit script applying the Seraph_Profits market-data refresh
# (scheduled runner updated currentAveragePrice* / Leve profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 228.71428
$ws.Range("I6").Value = 237.5
$ws.Range("J6").Value = 217
$ws.Range("K6").Value = 712.5
$ws.Range("L6").Value = 651
$ws.Range("M6").Value = -600.5
$ws.Range("N6").Value = -875

$ws.Range("H40").Value = 2847.3704
$ws.Range("I40").Value = 3313.8462
$ws.Range("J40").Value = 2414.2144
$ws.Range("K40").Value = 3313.8462
$ws.Range("L40").Value = 2414.2144
$ws.Range("M40").Value = -3138.8462
$ws.Range("N40").Value = -2764.2144

$ws.Range("H62").Value = 5788.8
$ws.Range("I62").Value = 3577.6
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3577.6
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -2953.6
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 5788.8
$ws.Range("I65").Value = 3577.6
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 17888
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -14768
$ws.Range("N65").Value = -46240

$ws.Range("H88").Value = 2069.3635
$ws.Range("I88").Value = 1862.3334
$ws.Range("J88").Value = 2147
$ws.Range("K88").Value = 1862.3334
$ws.Range("L88").Value = 2147
$ws.Range("M88").Value = -1456.3334
$ws.Range("N88").Value = -2959

$ws.Range("H91").Value = 2069.3635
$ws.Range("I91").Value = 1862.3334
$ws.Range("J91").Value = 2147
$ws.Range("K91").Value = 1862.3334
$ws.Range("L91").Value = 2147
$ws.Range("M91").Value = -458.3334
$ws.Range("N91").Value = -4955

$ws.Range("H96").Value = 1146.4286
$ws.Range("I96").Value = 1146.4286
$ws.Range("K96").Value = 3439.2858
$ws.Range("M96").Value = -2066.2858

$ws.Range("H99").Value = 341.6
$ws.Range("I99").Value = 302
$ws.Range("K99").Value = 906
$ws.Range("M99").Value = 592

$ws.Range("H101").Value = 430.75
$ws.Range("I101").Value = 430.75
$ws.Range("K101").Value = 1292.25
$ws.Range("M101").Value = 329.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2352.7693
$ws.Range("I61").Value = 2352.7693
$ws.Range("K61").Value = 2352.7693
$ws.Range("M61").Value = -2140.7693

$ws.Range("H102").Value = 1256.091
$ws.Range("I102").Value = 765.875
$ws.Range("K102").Value = 765.875
$ws.Range("M102").Value = 856.125

$ws.Range("H132").Value = 843.0625
$ws.Range("I132").Value = 843.0625
$ws.Range("K132").Value = 2529.1875
$ws.Range("M132").Value = 0.8125

$ws.Range("H136").Value = 2352.7693
$ws.Range("I136").Value = 2352.7693
$ws.Range("K136").Value = 7058.3079
$ws.Range("M136").Value = -4508.3079

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2035.5454
$ws.Range("I94").Value = 1939.1
$ws.Range("K94").Value = 1939.1
$ws.Range("M94").Value = -1488.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 64999
$ws.Range("J59").Value = 64999
$ws.Range("L59").Value = 64999
$ws.Range("N59").Value = -67289

$ws.Range("H94").Value = 704.8
$ws.Range("I94").Value = 681
$ws.Range("K94").Value = 681
$ws.Range("M94").Value = -230

$ws.Range("H132").Value = 1515.5
$ws.Range("I132").Value = 1515.5
$ws.Range("K132").Value = 4546.5
$ws.Range("M132").Value = -2016.5

$ws.Range("H134").Value = 3260.3684
$ws.Range("I134").Value = 3222.818
$ws.Range("K134").Value = 9668.454000000002
$ws.Range("M134").Value = -7133.454000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2125508.5
$ws.Range("J4").Value = 669
$ws.Range("L4").Value = 2007
$ws.Range("N4").Value = -2231

$ws.Range("H7").Value = 20000490
$ws.Range("I7").Value = 50000500
$ws.Range("J7").Value = 484
$ws.Range("K7").Value = 150001500
$ws.Range("L7").Value = 1452
$ws.Range("M7").Value = -150001388
$ws.Range("N7").Value = -1676

$ws.Range("H34").Value = 964.7273
$ws.Range("I34").Value = 401.33334
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 1204.00002
$ws.Range("L34").Value = 10500
$ws.Range("M34").Value = -1120.00002
$ws.Range("N34").Value = -10668

$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -5706
$ws.Range("N39").ClearContents()

$ws.Range("H51").Value = 2093.3333
$ws.Range("I51").Value = 1850
$ws.Range("J51").Value = 2336.6667
$ws.Range("K51").Value = 5550
$ws.Range("L51").Value = 7010.000100000001
$ws.Range("M51").Value = -5090
$ws.Range("N51").Value = -7930.000100000001

$ws.Range("H55").Value = 4172.1113
$ws.Range("J55").Value = 4172.1113
$ws.Range("L55").Value = 12516.3339
$ws.Range("N55").Value = -12870.3339

$ws.Range("H88").Value = 4838
$ws.Range("I88").Value = 4838
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 14514
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -14086
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 4838
$ws.Range("I91").Value = 4838
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 14514
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -13032
$ws.Range("N91").ClearContents()

$ws.Range("H136").Value = 20000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H139").Value = 1694
$ws.Range("I139").Value = 1694
$ws.Range("K139").Value = 5082
$ws.Range("M139").Value = 58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2774.1428
$ws.Range("I102").Value = 2027.1333
$ws.Range("K102").Value = 2027.1333
$ws.Range("M102").Value = -405.1333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 22000
$ws.Range("J95").Value = 22000
$ws.Range("L95").Value = 22000
$ws.Range("N95").Value = -27492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 55649.75
$ws.Range("J64").Value = 55649.75
$ws.Range("L64").Value = 55649.75
$ws.Range("N64").Value = -56145.75

$ws.Range("H67").Value = 55649.75
$ws.Range("J67").Value = 55649.75
$ws.Range("L67").Value = 55649.75
$ws.Range("N67").Value = -57365.75

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 2460.375
$ws.Range("I122").Value = 2383.2856
$ws.Range("K122").Value = 7149.8568
$ws.Range("M122").Value = -4699.8568

$ws.Range("H132").Value = 2874.244
$ws.Range("I132").Value = 2395.1562
$ws.Range("J132").Value = 4577.6665
$ws.Range("K132").Value = 7185.4686
$ws.Range("L132").Value = 13732.9995
$ws.Range("M132").Value = -4655.4686
$ws.Range("N132").Value = -18792.9995

$ws.Range("H136").Value = 1632.7778
$ws.Range("I136").Value = 992.1429000000001
$ws.Range("K136").Value = 2976.4287
$ws.Range("M136").Value = -426.4287000000004
